$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 131137298
$ws.Range("B4").Value = 57064
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 102612
$ws.Range("F4").Value = 'Järpe'
$ws.Range("G4").Value = 'Tetrastes bonasia'
$ws.Range("H4").Value = '(Linnaeus, 1758)'
$ws.Range("P4").Value = 'Kvarnfallet, Dlr'
$ws.Range("Q4").Value = 480012
$ws.Range("R4").Value = 6674519
$ws.Range("S4").Value = 5
$ws.Range("T4").Value = 'Dalarna'
$ws.Range("U4").Value = 'Ludvika'
$ws.Range("V4").Value = 'Dalarna'
$ws.Range("W4").Value = 'Grangärde'
$ws.Range("Y4").Formula = "=TEXT(DATE(2026,2,13),""yyyy-mm-dd"")"
$ws.Range("Z4").Value = '10:57'
$ws.Range("AA4").Formula = "=TEXT(DATE(2026,2,13),""yyyy-mm-dd"")"
$ws.Range("AB4").Value = '10:57'
$ws.Range("AC4").Value = '1 par födosökande'
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = 'Tobias Hellberg'
$ws.Range("AX4").Value = 'Tobias Hellberg'

# Row 5
$ws.Range("A5").Value = 131137846
$ws.Range("B5").Value = 57884
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = 'Tretåig hackspett'
$ws.Range("G5").Value = 'Picoides tridactylus'
$ws.Range("H5").Value = '(Linnaeus, 1758)'
$ws.Range("L5").Value = 'hona'
$ws.Range("M5").Value = 'födosökande'
$ws.Range("N5").Value = 'observerad'
$ws.Range("P5").Value = 'Kvarnfallet, Dlr'
$ws.Range("Q5").Value = 480066
$ws.Range("R5").Value = 6674279
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = 'Dalarna'
$ws.Range("U5").Value = 'Ludvika'
$ws.Range("V5").Value = 'Dalarna'
$ws.Range("W5").Value = 'Grangärde'
$ws.Range("Y5").Formula = "=TEXT(DATE(2026,2,13),""yyyy-mm-dd"")"
$ws.Range("Z5").Value = '11:56'
$ws.Range("AA5").Formula = "=TEXT(DATE(2026,2,13),""yyyy-mm-dd"")"
$ws.Range("AB5").Value = '11:56'
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = 'Tobias Hellberg'
$ws.Range("AX5").Value = 'Tobias Hellberg'

# Row 6
$ws.Range("A6").Value = 131137572
$ws.Range("B6").Value = 58043
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 103021
$ws.Range("F6").Value = 'Talltita'
$ws.Range("G6").Value = 'Poecile montanus'
$ws.Range("H6").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("N6").Value = 'observerad'
$ws.Range("P6").Value = 'Kvarnfallet, Dlr'
$ws.Range("Q6").Value = 480040
$ws.Range("R6").Value = 6674318
$ws.Range("S6").Value = 5
$ws.Range("T6").Value = 'Dalarna'
$ws.Range("U6").Value = 'Ludvika'
$ws.Range("V6").Value = 'Dalarna'
$ws.Range("W6").Value = 'Grangärde'
$ws.Range("Y6").Formula = "=TEXT(DATE(2026,2,13),""yyyy-mm-dd"")"
$ws.Range("Z6").Value = '11:33'
$ws.Range("AA6").Formula = "=TEXT(DATE(2026,2,13),""yyyy-mm-dd"")"
$ws.Range("AB6").Value = '11:33'
$ws.Range("AC6").Value = '2 st'
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = 'Tobias Hellberg'
$ws.Range("AX6").Value = 'Tobias Hellberg'

# Row 7
$ws.Range("A7").Value = 131137857
$ws.Range("B7").Value = 58043
$ws.Range("D7").Value = 'NT'
$ws.Range("E7").Value = 103021
$ws.Range("F7").Value = 'Talltita'
$ws.Range("G7").Value = 'Poecile montanus'
$ws.Range("H7").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("N7").Value = 'observerad'
$ws.Range("P7").Value = 'Kvarnfallet, Dlr'
$ws.Range("Q7").Value = 480066
$ws.Range("R7").Value = 6674279
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = 'Dalarna'
$ws.Range("U7").Value = 'Ludvika'
$ws.Range("V7").Value = 'Dalarna'
$ws.Range("W7").Value = 'Grangärde'
$ws.Range("Y7").Formula = "=TEXT(DATE(2026,2,13),""yyyy-mm-dd"")"
$ws.Range("Z7").Value = '11:59'
$ws.Range("AA7").Formula = "=TEXT(DATE(2026,2,13),""yyyy-mm-dd"")"
$ws.Range("AB7").Value = '11:59'
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = 'Tobias Hellberg'
$ws.Range("AX7").Value = 'Tobias Hellberg'

# Row 8
$ws.Range("A8").Value = 131137098
$ws.Range("B8").Value = 80348
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 6458
$ws.Range("F8").Value = 'Lunglav'
$ws.Range("G8").Value = 'Lobaria pulmonaria'
$ws.Range("H8").Value = '(L.) Hoffm.'
$ws.Range("P8").Value = 'Kvarnfallet, Dlr'
$ws.Range("Q8").Value = 479869
$ws.Range("R8").Value = 6674500
$ws.Range("S8").Value = 5
$ws.Range("T8").Value = 'Dalarna'
$ws.Range("U8").Value = 'Ludvika'
$ws.Range("V8").Value = 'Dalarna'
$ws.Range("W8").Value = 'Grangärde'
$ws.Range("Y8").Formula = "=TEXT(DATE(2026,2,13),""yyyy-mm-dd"")"
$ws.Range("Z8").Value = '10:26'
$ws.Range("AA8").Formula = "=TEXT(DATE(2026,2,13),""yyyy-mm-dd"")"
$ws.Range("AB8").Value = '10:26'
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AJ8").Value = 'sälg'
$ws.Range("AK8").Value = 'Salix caprea'
$ws.Range("AO8").Value = 'Salix caprea'
$ws.Range("AW8").Value = 'Tobias Hellberg'
$ws.Range("AX8").Value = 'Tobias Hellberg'

# Freeze the TEXT() formulas above into literal text values (matches source export,
# which stores dates as plain text rather than formulas or serials).
$ws.Range("Y4").Copy()
$ws.Range("Y4").PasteSpecial(-4163)
$ws.Range("AA4").Copy()
$ws.Range("AA4").PasteSpecial(-4163)
$ws.Range("Y5").Copy()
$ws.Range("Y5").PasteSpecial(-4163)
$ws.Range("AA5").Copy()
$ws.Range("AA5").PasteSpecial(-4163)
$ws.Range("Y6").Copy()
$ws.Range("Y6").PasteSpecial(-4163)
$ws.Range("AA6").Copy()
$ws.Range("AA6").PasteSpecial(-4163)
$ws.Range("Y7").Copy()
$ws.Range("Y7").PasteSpecial(-4163)
$ws.Range("AA7").Copy()
$ws.Range("AA7").PasteSpecial(-4163)
$ws.Range("Y8").Copy()
$ws.Range("Y8").PasteSpecial(-4163)
$ws.Range("AA8").Copy()
$ws.Range("AA8").PasteSpecial(-4163)
$excel.CutCopyMode = $false

Write-Host "Added rows 4-8 to Artfynd sheet."
